# Replace the single run
#     "Add no-image place holder image"
# with four separate runs that together read
#     "Check  no-image place holder and no title"
# i.e.
#     "Check " / " no-image place holder" / " and" / " no title"
#
# A plain Find/Replace (or Range.Text = "...") would produce one merged run,
# and the diff we need to reproduce keeps the text split across four <w:r>
# elements (no run formatting differs between them - Word just happens to
# have typed/pasted this in pieces). To get distinct <w:r> elements out of
# this engine's serializer (which always coalesces adjacent runs that carry
# identical formatting) we briefly turn the text into separate paragraphs -
# each chunk gets its own paragraph mark, so each chunk is necessarily its
# own run - and then delete those temporary paragraph marks again to merge
# everything back into the original paragraph. Deleting a paragraph mark
# merges it into the *next* paragraph and keeps that next paragraph's mark
# (and therefore its identity: w14:paraId, rsidR, ...), while leaving the
# text that came from each side as distinct runs. So we build the new
# paragraphs *before* the original one (new chunks first, original
# paragraph/mark last) so that after all the merges collapse back down, the
# surviving paragraph is still the original one.

$d = $word.ActiveDocument

# Locate the exact text to replace and remember its Range bounds.
$target = $d.Content
$null = $target.Find.Execute("Add no-image place holder image")
$origStart = $target.Start
$origEnd = $target.End

# 1-based index of the paragraph holding that text (paragraph boundaries
# before $origStart are unaffected by anything we do below, so this stays
# valid throughout).
$paraIndex = $d.Range(0, $origStart).Paragraphs.Count + 1

# The runs the new text is split across, in reading order.
$segments = @("Check ", " no-image place holder", " and", " no title")

# Insert every segment but the last as its own new paragraph, right before
# the original text, each terminated with a paragraph mark.
$cursor = $origStart
for ($i = 0; $i -lt $segments.Length - 1; $i++) {
    $chunk = $segments[$i] + "`r"
    $insPoint = $d.Range($cursor, $cursor)
    $insPoint.InsertBefore($chunk)
    $cursor = $cursor + $chunk.Length
}

# Replace the (now shifted-forward) original text with the final segment.
$shift = $cursor - $origStart
$remaining = $d.Range($cursor, $origEnd + $shift)
$remaining.Text = $segments[$segments.Length - 1]

# Merge the temporary paragraph breaks back out, restoring a single
# paragraph that still carries the original paragraph's identity, with the
# four text chunks left behind as four separate runs.
for ($i = 0; $i -lt $segments.Length - 1; $i++) {
    $para = $d.Paragraphs.Item($paraIndex)
    $pilcrowPos = $para.Range.End - 1
    $d.Range($pilcrowPos, $pilcrowPos + 1).Delete()
}
